$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "G2"
$ws.Range("B4").Value = "Mask1"
$ws.Range("C4").Value = 45860
$ws.Range("C4").NumberFormat = "YYYY-MM-DD"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
